# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for the
# rows whose dialog-act classification changed after the transcript clean up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(4, 9).Value = 'sd'
$ws.Cells.Item(4, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(5, 9).Value = 'b'
$ws.Cells.Item(5, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(9, 9).Value = 'sd'
$ws.Cells.Item(9, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(14, 9).Value = 'sd'
$ws.Cells.Item(14, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(19, 9).Value = 'b'
$ws.Cells.Item(19, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(30, 9).Value = 'aa'
$ws.Cells.Item(30, 10).Value = 'Agree/Accept'
$ws.Cells.Item(53, 9).Value = 'aa'
$ws.Cells.Item(53, 10).Value = 'Agree/Accept'
$ws.Cells.Item(91, 9).Value = 'b'
$ws.Cells.Item(91, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(99, 9).Value = 'sv'
$ws.Cells.Item(99, 10).Value = 'Statement-opinion'
$ws.Cells.Item(101, 9).Value = 'b'
$ws.Cells.Item(101, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(102, 9).Value = 'ba'
$ws.Cells.Item(102, 10).Value = 'Appreciation'
$ws.Cells.Item(107, 9).Value = 'sd'
$ws.Cells.Item(107, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(109, 9).Value = 'b'
$ws.Cells.Item(109, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(112, 9).Value = 'b'
$ws.Cells.Item(112, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(131, 9).Value = 'b'
$ws.Cells.Item(131, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(134, 9).Value = 'sd'
$ws.Cells.Item(134, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(135, 9).Value = 'sd'
$ws.Cells.Item(135, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(146, 9).Value = 'aa'
$ws.Cells.Item(146, 10).Value = 'Agree/Accept'
$ws.Cells.Item(159, 9).Value = 'sv'
$ws.Cells.Item(159, 10).Value = 'Statement-opinion'
$ws.Cells.Item(171, 9).Value = '%'
$ws.Cells.Item(171, 10).Value = 'Uninterpretable'
$ws.Cells.Item(184, 9).Value = 'sv'
$ws.Cells.Item(184, 10).Value = 'Statement-opinion'
$ws.Cells.Item(186, 9).Value = 'sd'
$ws.Cells.Item(186, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(188, 9).Value = 'aa'
$ws.Cells.Item(188, 10).Value = 'Agree/Accept'
$ws.Cells.Item(189, 9).Value = 'sv'
$ws.Cells.Item(189, 10).Value = 'Statement-opinion'
$ws.Cells.Item(198, 9).Value = 'sv'
$ws.Cells.Item(198, 10).Value = 'Statement-opinion'
$ws.Cells.Item(204, 9).Value = '%'
$ws.Cells.Item(204, 10).Value = 'Uninterpretable'
$ws.Cells.Item(205, 9).Value = 'sv'
$ws.Cells.Item(205, 10).Value = 'Statement-opinion'
$ws.Cells.Item(206, 9).Value = 'sd'
$ws.Cells.Item(206, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(242, 9).Value = 'sd'
$ws.Cells.Item(242, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(243, 9).Value = 'sd'
$ws.Cells.Item(243, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(251, 9).Value = 'sd'
$ws.Cells.Item(251, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(257, 9).Value = 'sd'
$ws.Cells.Item(257, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(281, 9).Value = 'sd'
$ws.Cells.Item(281, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(283, 9).Value = 'sd'
$ws.Cells.Item(283, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(285, 9).Value = 'sd'
$ws.Cells.Item(285, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(287, 9).Value = 'aa'
$ws.Cells.Item(287, 10).Value = 'Agree/Accept'
$ws.Cells.Item(294, 9).Value = 'aa'
$ws.Cells.Item(294, 10).Value = 'Agree/Accept'
$ws.Cells.Item(315, 9).Value = '%'
$ws.Cells.Item(315, 10).Value = 'Uninterpretable'
$ws.Cells.Item(327, 9).Value = 'sd'
$ws.Cells.Item(327, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(339, 9).Value = 'sd'
$ws.Cells.Item(339, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(352, 9).Value = 'sd'
$ws.Cells.Item(352, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(353, 9).Value = 'sv'
$ws.Cells.Item(353, 10).Value = 'Statement-opinion'
$ws.Cells.Item(354, 9).Value = 'sv'
$ws.Cells.Item(354, 10).Value = 'Statement-opinion'
$ws.Cells.Item(362, 9).Value = 'sd'
$ws.Cells.Item(362, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(363, 9).Value = 'sv'
$ws.Cells.Item(363, 10).Value = 'Statement-opinion'
$ws.Cells.Item(364, 9).Value = 'sv'
$ws.Cells.Item(364, 10).Value = 'Statement-opinion'
$ws.Cells.Item(367, 9).Value = 'sv'
$ws.Cells.Item(367, 10).Value = 'Statement-opinion'
$ws.Cells.Item(369, 9).Value = 'sv'
$ws.Cells.Item(369, 10).Value = 'Statement-opinion'
